$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2595
$ws.Range("I62").Value = 2397.5
$ws.Range("J62").Value = 2990
$ws.Range("K62").Value = 2397.5
$ws.Range("L62").Value = 2990
$ws.Range("M62").Value = -1773.5
$ws.Range("N62").Value = -4238
$ws.Range("H65").Value = 2595
$ws.Range("I65").Value = 2397.5
$ws.Range("J65").Value = 2990
$ws.Range("K65").Value = 11987.5
$ws.Range("L65").Value = 14950
$ws.Range("M65").Value = -8867.5
$ws.Range("N65").Value = -21190
$ws.Range("H80").Value = 1163.9166
$ws.Range("I80").Value = 1602.7142
$ws.Range("J80").Value = 549.6
$ws.Range("K80").Value = 4808.142599999999
$ws.Range("L80").Value = 1648.8
$ws.Range("M80").Value = -3810.142599999999
$ws.Range("N80").Value = -3644.8
$ws.Range("H83").Value = 1163.9166
$ws.Range("I83").Value = 1602.7142
$ws.Range("J83").Value = 549.6
$ws.Range("K83").Value = 14424.4278
$ws.Range("L83").Value = 4946.400000000001
$ws.Range("M83").Value = -9432.427799999999
$ws.Range("N83").Value = -14930.4
$ws.Range("H86").Value = 2557
$ws.Range("I86").Value = 2399
$ws.Range("K86").Value = 2399
$ws.Range("M86").Value = -1276
$ws.Range("H89").Value = 2557
$ws.Range("I89").Value = 2399
$ws.Range("K89").Value = 11995
$ws.Range("M89").Value = -6379
$ws.Range("H131").Value = 2063.5293
$ws.Range("J131").Value = 4015.7144
$ws.Range("L131").Value = 12047.1432
$ws.Range("N131").Value = -22127.1432
$ws.Range("H138").Value = 1875.2174
$ws.Range("J138").Value = 2169.25
$ws.Range("L138").Value = 6507.75
$ws.Range("N138").Value = -16787.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 17900
$ws.Range("I28").Value = 6000
$ws.Range("J28").Value = 29800
$ws.Range("K28").Value = 6000
$ws.Range("L28").Value = 29800
$ws.Range("M28").Value = -5808
$ws.Range("N28").Value = -30184
$ws.Range("H97").Value = 503
$ws.Range("I97").Value = 509.625
$ws.Range("J97").Value = 450
$ws.Range("K97").Value = 509.625
$ws.Range("L97").Value = 450
$ws.Range("M97").Value = -13.625
$ws.Range("N97").Value = -1442
$ws.Range("H99").Value = 17900
$ws.Range("I99").Value = 6000
$ws.Range("J99").Value = 29800
$ws.Range("K99").Value = 6000
$ws.Range("L99").Value = 29800
$ws.Range("M99").Value = -3005
$ws.Range("N99").Value = -35790
$ws.Range("H132").Value = 4002
$ws.Range("I132").Value = 2506
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 7518
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -4988
$ws.Range("N132").Value = -20058.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3031.4666
$ws.Range("I107").Value = 2028.1428
$ws.Range("J107").Value = 3909.375
$ws.Range("K107").Value = 2028.1428
$ws.Range("L107").Value = 3909.375
$ws.Range("M107").Value = -108.1428000000001
$ws.Range("N107").Value = -7749.375
$ws.Range("H134").Value = 11501.333
$ws.Range("I134").Value = 12655.357
$ws.Range("K134").Value = 37966.071
$ws.Range("M134").Value = -35431.071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2367.366
$ws.Range("I31").Value = 1130.8518
$ws.Range("J31").Value = 4752.0713
$ws.Range("K31").Value = 1130.8518
$ws.Range("L31").Value = 4752.0713
$ws.Range("M31").Value = -835.8517999999999
$ws.Range("N31").Value = -5342.0713
$ws.Range("H34").Value = 2367.366
$ws.Range("I34").Value = 1130.8518
$ws.Range("J34").Value = 4752.0713
$ws.Range("K34").Value = 1130.8518
$ws.Range("L34").Value = 4752.0713
$ws.Range("M34").Value = -928.8517999999999
$ws.Range("N34").Value = -5156.0713
$ws.Range("H58").Value = 1282.4615
$ws.Range("I58").Value = 1211.75
$ws.Range("K58").Value = 1211.75
$ws.Range("M58").Value = -1008.75
$ws.Range("H86").Value = 1325.6
$ws.Range("I86").Value = 1325.6
$ws.Range("K86").Value = 1325.6
$ws.Range("M86").Value = -202.5999999999999
$ws.Range("H89").Value = 1325.6
$ws.Range("I89").Value = 1325.6
$ws.Range("K89").Value = 6628
$ws.Range("M89").Value = -1012
$ws.Range("H134").Value = 998.2308
$ws.Range("I134").Value = 951.2727
$ws.Range("K134").Value = 2853.8181
$ws.Range("M134").Value = -318.8181
$ws.Range("H136").Value = 1282.4615
$ws.Range("I136").Value = 1211.75
$ws.Range("K136").Value = 3635.25
$ws.Range("M136").Value = -1085.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 748.25
$ws.Range("J98").Value = 796.5
$ws.Range("L98").Value = 2389.5
$ws.Range("N98").Value = -5385.5
$ws.Range("H131").Value = 12215299
$ws.Range("J131").Value = 23549.143
$ws.Range("L131").Value = 70647.429
$ws.Range("N131").Value = -80727.429
$ws.Range("H137").Value = 3743.9443
$ws.Range("I137").Value = 1186.9
$ws.Range("J137").Value = 6940.25
$ws.Range("K137").Value = 3560.7
$ws.Range("L137").Value = 20820.75
$ws.Range("M137").Value = 1539.3
$ws.Range("N137").Value = -31020.75
$ws.Range("H139").Value = 8267.799999999999
$ws.Range("I139").Value = 9155.308000000001
$ws.Range("J139").Value = 2499
$ws.Range("K139").Value = 27465.924
$ws.Range("L139").Value = 7497
$ws.Range("M139").Value = -22325.924
$ws.Range("N139").Value = -17777

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2997.5
$ws.Range("I80").Value = 2995
$ws.Range("K80").Value = 2995
$ws.Range("M80").Value = -1997
$ws.Range("H83").Value = 2997.5
$ws.Range("I83").Value = 2995
$ws.Range("K83").Value = 14975
$ws.Range("M83").Value = -9983
$ws.Range("H132").Value = 3992.1667
$ws.Range("I132").Value = 3266
$ws.Range("K132").Value = 9798
$ws.Range("M132").Value = -7268

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6322.077
$ws.Range("I7").Value = 3260.2
$ws.Range("J7").Value = 8235.75
$ws.Range("K7").Value = 3260.2
$ws.Range("L7").Value = 8235.75
$ws.Range("M7").Value = -3148.2
$ws.Range("N7").Value = -8459.75
$ws.Range("H55").Value = 382
$ws.Range("I55").Value = 120
$ws.Range("J55").Value = 480.25
$ws.Range("K55").Value = 120
$ws.Range("L55").Value = 480.25
$ws.Range("M55").Value = 53
$ws.Range("N55").Value = -826.25
$ws.Range("H68").Value = 3496.6667
$ws.Range("J68").Value = 5000
$ws.Range("L68").Value = 5000
$ws.Range("N68").Value = -6498
$ws.Range("H71").Value = 3496.6667
$ws.Range("J71").Value = 5000
$ws.Range("L71").Value = 25000
$ws.Range("N71").Value = -32488
$ws.Range("H126").Value = 6322.077
$ws.Range("I126").Value = 3260.2
$ws.Range("J126").Value = 8235.75
$ws.Range("K126").Value = 9780.599999999999
$ws.Range("L126").Value = 24707.25
$ws.Range("M126").Value = -7310.599999999999
$ws.Range("N126").Value = -29647.25
$ws.Range("H132").Value = 1316.25
$ws.Range("I132").Value = 1023.9286
$ws.Range("J132").Value = 1608.5714
$ws.Range("K132").Value = 3071.7858
$ws.Range("L132").Value = 4825.7142
$ws.Range("M132").Value = -541.7857999999997
$ws.Range("N132").Value = -9885.7142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H63").Value = 1677075.4
$ws.Range("J63").Value = 2512500
$ws.Range("L63").Value = 2512500
$ws.Range("N63").Value = -2513748
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H66").Value = 1677075.4
$ws.Range("J66").Value = 2512500
$ws.Range("L66").Value = 7537500
$ws.Range("N66").Value = -7543740
$ws.Range("H132").Value = 6599.75
$ws.Range("I132").Value = 299
$ws.Range("J132").Value = 7499.857
$ws.Range("K132").Value = 897
$ws.Range("L132").Value = 22499.571
$ws.Range("M132").Value = 1633
$ws.Range("N132").Value = -27559.571
